$d = $word.ActiveDocument

# 1) Remove the trailing " Report #01" run from the title block.
$d.Content.Find.Execute(" Report #01", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2) Split ", Terence Henriod" into ", Ter" / "ence Henriod" and move the
#    "_GoBack" bookmark to the split point.
$r2 = $d.Content
$found2 = $r2.Find.Execute(", Ter")
if ($found2) {
    $splitPoint = $d.Range($r2.End, $r2.End)
    $d.Bookmarks.Add("_GoBack", $splitPoint)
}

# 3) Merge the two runs around the old "_GoBack" bookmark location back
#    into a single run (removing the bookmark from there).
$d.Content.Find.Execute("before the competition/demonstration.", $true, $false, $false, $false, $false, $true, 1, $false, "before the competition/demonstration.", 2) | Out-Null
